$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(43526, 43527, 43528, 43529, 43530, 43531, 43532)

$row = 144
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = "Recife"
    $ws.Cells.Item($row, 2).Value = $d
    $ws.Cells.Item($row, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 4).Value = 983
    $ws.Cells.Item($row, 5).Value = 8
    $row = $row + 1
}

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E150"))

$ws.Range("B149").Select()
